# Refresh chart data with updated forecast values on a timeframe/schedule basis.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B (forecast) and column C (hour) per row.
$updates = @{
    2  = @{ B = 1281.8; C = 108 }
    3  = @{ B = 523;    C = 108 }
    4  = @{ B = 19;     C = 108 }
    5  = @{ B = 191;    C = 108 }
    6  = @{ B = 96;     C = 108 }
    7  = @{ B = 193;    C = 108 }
    10 = @{ B = 957;    C = 108 }
    11 = @{ B = 321;    C = 108 }
    12 = @{ B = 1137.5; C = 108 }
    13 = @{ B = 816;    C = 108 }
    14 = @{ B = 1148;   C = 108 }
    15 = @{ B = 371;    C = 108 }
    17 = @{ B = 202;    C = 108 }
    18 = @{ B = 6;      C = 108 }
    19 = @{ B = 18;     C = 108 }
    20 = @{ B = 139;    C = 108 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
}

$wb.Save()
